$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# --- Row 12: reuse formatting from row 11 (A=style4, B=default, C=style5, D=style4) ---
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null

# --- Row 13 & 14: column A/C/D keep the row-11 pattern, column B additionally styled (style4) ---
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A13:A14").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("B13:B14").PasteSpecial(-4122) | Out-Null

$ws.Range("C11").Copy() | Out-Null
$ws.Range("C13:C14").PasteSpecial(-4122) | Out-Null

$ws.Range("D11").Copy() | Out-Null
$ws.Range("D13:D14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Populate the new Asset rows: LeadSheetURL, ReportFileID, ReportFileURL ---
$ws.Range("D12").Value = "Gdrive file URL for the Lead Spreadsheet document"
$ws.Range("B12").Value = "13_LeadSheetURL"
$ws.Range("A12").Value = "LeadSheetURL"
$ws.Range("C12").Value = "Shared"

$ws.Range("A13").Value = "ReportFileID"
$ws.Range("A14").Value = "ReportFileURL"
$ws.Range("B13").Value = "13_ReportFileID"
$ws.Range("B14").Value = "13_ReportFileURL"
$ws.Range("C13").Value = "Shared"
$ws.Range("C14").Value = "Shared"
$ws.Range("D13").Value = "Gdrive file ID for the current execution report"
$ws.Range("D14").Value = "Gdrive file URL for the current execution report"

# --- Move the active selection the way the author left it ---
$ws.Range("D17").Select() | Out-Null
